# Apply the "gap analysis" partial commit:
#  - insert a new "GapAnalysis" worksheet between "Test Data" and "Alignment Type"
#  - populate it with the gap-analysis summary table
#  - leave the other two sheets' data untouched (only their tab/position shift)

$wb = $excel.ActiveWorkbook
$wsTestData = $wb.Worksheets.Item(1)
$wsAlignType = $wb.Worksheets.Item(2)

# Insert the new sheet right after "Test Data" (i.e. before "Alignment Type")
$ws = $wb.Worksheets.Add($null, $wsTestData)
$ws.Name = "GapAnalysis"

# --- populate cells (order mirrors how the source workbook was authored) ---
$ws.Range("A3").Value = "Automation Report Test-2"
$ws.Range("A2").Value = "Parent"
$ws.Range("B3").Value = "Reading"
$ws.Range("C3").Value = "Automation Report Test (US Grade K-12)"
$ws.Range("D3").Value = "TEST 1 (US Grade K-12)"
$ws.Range("B2").Value = "Child-Mapping"
$ws.Range("E3").Value = "Exact"
$ws.Range("F2").Value = "Add Note(Comments)"
$ws.Range("F3").Value = " Gap Analysis "
$ws.Range("B4").Value = "Text Types and Purposes - LEVEL 2."
$ws.Range("A1").Value = "SOURCE"
$ws.Range("C1").Value = "TARGET"
$ws.Range("D4").Value = "TEST 2 (US Grade K-12)"
$ws.Range("E4").Value = "Related"
$ws.Range("B5").Value = "Write arguments to support claims in an analysis of substantive topics or texts, using valid reasoning and relevant and sufficient evidence - LEVEL 3."
$ws.Range("D5").Value = "TEST 3 (US Grade K-12)"
$ws.Range("E5").Value = "Broad"
$ws.Range("B6").Value = "Reading LEVEL 4."
$ws.Range("E6").Value = "Close"
$ws.Range("D6").Value = "TEST 4 (US Grade K-12)"
$ws.Range("B7").Value = "Reading LOWEST LEVEL"
$ws.Range("D7").Value = "Conventions of Standard English (US Grade K-12)"
$ws.Range("E7").Value = "Narrow"
$ws.Range("F7").Value = "Lowest node alignment"

# C2/D2 repeat the row-2 labels (reuse the same shared strings as A2/B2)
$ws.Range("C2").Value = "Parent"
$ws.Range("D2").Value = "Child-Mapping"
# E2 reuses the pre-existing "Alignment Type" shared string
$ws.Range("E2").Value = "Alignment Type"

# --- formatting -------------------------------------------------------
$headerCells = @("A1","C1","A2","B2","C2","D2","E2","F2")
foreach ($addr in $headerCells) {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.WrapText = $true
}

$dataCells = @(
    "A3","B3","C3","D3","E3","F3",
    "B4","D4","E4",
    "B5","D5","E5",
    "B6","D6","E6",
    "B7","D7","E7","F7"
)
foreach ($addr in $dataCells) {
    $ws.Range($addr).WrapText = $true
}

# Row heights for the two wrapped long-text rows
$ws.Rows.Item(5).RowHeight = 75
$ws.Rows.Item(7).RowHeight = 30

# Column widths (autofit-like widths as saved in the source file)
$ws.Columns.Item(1).ColumnWidth = 23.43
$ws.Columns.Item(2).ColumnWidth = 35.43
$ws.Columns.Item(3).ColumnWidth = 36.28
$ws.Columns.Item(4).ColumnWidth = 36.28
$ws.Columns.Item(5).ColumnWidth = 19.94
$ws.Columns.Item(6).ColumnWidth = 19.8
$ws.Columns.Item(7).ColumnWidth = 19.8

# Selection left on D7, as in the saved file; GapAnalysis ends up the active tab
$ws.Range("D7").Select() | Out-Null
$ws.Activate() | Out-Null
